# Auto-generated edit script: apply scheduled market-price refresh to Asura_Profits workbook
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H-N) for specific leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 785.5
$ws.Range("I19").Value = 633.3333
$ws.Range("J19").Value = 876.8
$ws.Range("K19").Value = 633.3333
$ws.Range("L19").Value = 876.8
$ws.Range("M19").Value = -458.3333
$ws.Range("N19").Value = -1226.8

$ws.Range("H80").Value = 7260.706
$ws.Range("I80").Value = 1042
$ws.Range("J80").Value = 12788.444
$ws.Range("K80").Value = 3126
$ws.Range("L80").Value = 38365.33199999999
$ws.Range("M80").Value = -2128
$ws.Range("N80").Value = -40361.33199999999

$ws.Range("H83").Value = 7260.706
$ws.Range("I83").Value = 1042
$ws.Range("J83").Value = 12788.444
$ws.Range("K83").Value = 9378
$ws.Range("L83").Value = 115095.996
$ws.Range("M83").Value = -4386
$ws.Range("N83").Value = -125079.996

$ws.Range("H86").Value = 2437.7778
$ws.Range("I86").Value = 2563
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 2563
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = -1440
$ws.Range("N86").Value = -4245.5

$ws.Range("H89").Value = 2437.7778
$ws.Range("I89").Value = 2563
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 12815
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = -7199
$ws.Range("N89").Value = -21229.5

$ws.Range("H137").Value = 1454.2174
$ws.Range("I137").Value = 1328.0714
$ws.Range("J137").Value = 1650.4445
$ws.Range("K137").Value = 3984.2142
$ws.Range("L137").Value = 4951.333500000001
$ws.Range("M137").Value = -1434.2142
$ws.Range("N137").Value = -10051.3335

$ws.Range("H138").Value = 2099.35
$ws.Range("I138").Value = 1192.8857
$ws.Range("J138").Value = 4214.433
$ws.Range("K138").Value = 3578.6571
$ws.Range("L138").Value = 12643.299
$ws.Range("M138").Value = 1561.3429
$ws.Range("N138").Value = -22923.299

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12792.167
$ws.Range("I32").Value = 13090.83
$ws.Range("J32").Value = 10786.857
$ws.Range("K32").Value = 13090.83
$ws.Range("L32").Value = 10786.857
$ws.Range("M32").Value = -12803.83
$ws.Range("N32").Value = -11360.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 32000
$ws.Range("J69").Value = 32000
$ws.Range("L69").Value = 32000
$ws.Range("N69").Value = -33622

$ws.Range("H72").Value = 32000
$ws.Range("J72").Value = 32000
$ws.Range("L72").Value = 96000
$ws.Range("N72").Value = -104112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2355.0833
$ws.Range("I31").Value = 1723.174
$ws.Range("J31").Value = 3473.077
$ws.Range("K31").Value = 1723.174
$ws.Range("L31").Value = 3473.077
$ws.Range("M31").Value = -1428.174
$ws.Range("N31").Value = -4063.077

$ws.Range("H34").Value = 2355.0833
$ws.Range("I34").Value = 1723.174
$ws.Range("J34").Value = 3473.077
$ws.Range("K34").Value = 1723.174
$ws.Range("L34").Value = 3473.077
$ws.Range("M34").Value = -1521.174
$ws.Range("N34").Value = -3877.077

$ws.Range("H134").Value = 1554.05
$ws.Range("I134").Value = 1339.4517
$ws.Range("J134").Value = 2293.2222
$ws.Range("K134").Value = 4018.3551
$ws.Range("L134").Value = 6879.6666
$ws.Range("M134").Value = -1483.3551
$ws.Range("N134").Value = -11949.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1145.0741
$ws.Range("I5").Value = 1143.6364
$ws.Range("J5").Value = 1151.4
$ws.Range("K5").Value = 3430.9092
$ws.Range("L5").Value = 3454.2
$ws.Range("M5").Value = -3318.9092
$ws.Range("N5").Value = -3678.2

$ws.Range("H50").Value = 2157.4666
$ws.Range("I50").Value = 118
$ws.Range("J50").Value = 3942
$ws.Range("K50").Value = 354
$ws.Range("L50").Value = 11826
$ws.Range("M50").Value = 127
$ws.Range("N50").Value = -12788

$ws.Range("H53").Value = 2157.4666
$ws.Range("I53").Value = 118
$ws.Range("J53").Value = 3942
$ws.Range("K53").Value = 354
$ws.Range("L53").Value = 11826
$ws.Range("M53").Value = 127
$ws.Range("N53").Value = -12788

$ws.Range("H122").Value = 618
$ws.Range("I122").Value = 478.44446
$ws.Range("J122").Value = 732.1818
$ws.Range("K122").Value = 4306.00014
$ws.Range("L122").Value = 6589.6362
$ws.Range("M122").Value = -1856.00014
$ws.Range("N122").Value = -11489.6362

$ws.Range("H131").Value = 10640943
$ws.Range("I131").Value = 22306
$ws.Range("J131").Value = 11237496
$ws.Range("K131").Value = 66918
$ws.Range("L131").Value = 33712488
$ws.Range("M131").Value = -61878
$ws.Range("N131").Value = -33722568

$ws.Range("H135").Value = 1145.0741
$ws.Range("I135").Value = 1143.6364
$ws.Range("J135").Value = 1151.4
$ws.Range("K135").Value = 10292.7276
$ws.Range("L135").Value = 10362.6
$ws.Range("M135").Value = -7757.7276
$ws.Range("N135").Value = -15432.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1758.7333
$ws.Range("I132").Value = 1164.8518
$ws.Range("J132").Value = 2649.5557
$ws.Range("K132").Value = 3494.5554
$ws.Range("L132").Value = 7948.6671
$ws.Range("M132").Value = -964.5553999999997
$ws.Range("N132").Value = -13008.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 600.8570999999999
$ws.Range("I55").Value = 233.33333
$ws.Range("J55").Value = 876.5
$ws.Range("K55").Value = 233.33333
$ws.Range("L55").Value = 876.5
$ws.Range("M55").Value = -60.33332999999999
$ws.Range("N55").Value = -1222.5

$ws.Range("H61").Value = 28112.875
$ws.Range("I61").Value = 28112.875
$ws.Range("K61").Value = 28112.875
$ws.Range("M61").Value = -27910.875

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H113").Value = 28112.875
$ws.Range("I113").Value = 28112.875
$ws.Range("K113").Value = 28112.875
$ws.Range("M113").Value = -25942.875

$ws.Range("H130").Value = 22938
$ws.Range("J130").Value = 22938
$ws.Range("L130").Value = 22938
$ws.Range("N130").Value = -32978

$ws.Range("H136").Value = 2161.9795
$ws.Range("I136").Value = 1793.909
$ws.Range("J136").Value = 5401
$ws.Range("K136").Value = 5381.727000000001
$ws.Range("L136").Value = 16203
$ws.Range("M136").Value = -2831.727000000001
$ws.Range("N136").Value = -21303

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 186566
$ws.Range("J119").Value = 186566
$ws.Range("L119").Value = 186566
$ws.Range("N119").Value = -196242
